$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row: rename/insert/shift columns
#    Before: B=id_unico C=Nome D=Centro de Custo E=Salário F=Custo Gympass G=Custo Total
#    After:  B=id_unico C=Nome D=Centro de Custo E=Salario F=Custo Unimed G=Custo Gympass
#            H=Custo GitHub I=Custo Google Workspace J=Custo Claude K=Custo Total
# ---------------------------------------------------------------------
$ws.Range("K1").Value = "Custo Total"
$ws.Range("G1").Value = "Custo Gympass"
$ws.Range("E1").Value = "Salario"
$ws.Range("F1").Value = "Custo Unimed"
$ws.Range("H1").Value = "Custo GitHub"
$ws.Range("I1").Value = "Custo Google Workspace"
$ws.Range("J1").Value = "Custo Claude"

# Copy the header formatting (bold, border, centered) from an existing
# styled header cell (B1) onto the newly-introduced header cells.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:K1").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2. Data rows 2-63: one row per colaborador, resorted by id_unico (B)
#    ascending. Columns: B C D E F G H I J K
# ---------------------------------------------------------------------
$dataBlock = @"
00a0d290-74eb-5793-8df7-98f4280eb295|Antônia Pires|R&D|9976.53102463056|444.99|90.0|254.63|297.13|150.0|11213.28102463056
06be0738-b15f-5d07-892b-d45971647c79|Danilo Brandão|R&D|9095.090151715545|855.8|90.0|0.0|297.13|10.0|10348.02015171554
0eab9fc6-2c73-5288-af73-b13876c87b62|Catarina Chaves|R&D|3490.64446620121|444.99|90.0|0.0|297.13|10.0|4332.76446620121
1202a201-05dd-5f81-bff7-762f946c273e|Ariane Teixeira|COGS|8271.84079484528|533.99|90.0|254.63|297.13|150.0|9597.590794845279
1539ddc4-5241-5e9d-9086-cc8aea7dcbb1|Augusto Nunes|R&D|4963.54282637831|560.7|90.0|254.63|297.13|150.0|6316.00282637831
17533122-8361-509c-9e67-8887db2e0001|Diana Azevedo|R&D|6993.13850131002|-328.28|90.0|0.0|297.13||7051.98850131002
177d2a03-dd33-5953-8650-1de95414461f|Adolfo Moreira|R&D|3384.960765580439|533.99|90.0|254.63|297.13|150.0|4710.710765580439
25f4b46b-119b-54b7-a5fa-03827e5ed16f|Carolina Salgado|R&D|3266.657357480083|364.76|164.0|0.0|297.13|10.0|4102.547357480083
282b242b-4e87-51f2-a5ff-d141bd8825c8|Alessandra Dias|R&D|5318.987724937523|444.99|90.0|0.0|297.13|150.0|6301.107724937523
2bc2241b-333d-5e89-9b7c-7c6524872c01|Camila Dantas|R&D|1280.596453813949|444.99|90.0|254.63|297.13||2367.346453813949
2c297031-afd1-540c-a226-32291c3b420d|Eduardo Portela|R&D|1728.629771639805|560.7|117.0|0.0|297.13|50.0|2753.459771639805
2f144455-ef5b-5852-8ccb-71bbb361384b|Bruno Castro|R&D|4610.394047610112|364.76|90.0|254.63|297.13||5616.914047610113
361c4b67-5ca5-59d1-998e-08ceca3ebc57|Bárbara Neves|R&D|4587.767886848967|364.76|117.0|0.0|297.13|150.0|5516.657886848967
39f92dc0-2ae9-596f-a160-35025eebbff1|Afonso Barros|G&A|8962.7297218415|444.99|90.0|141.46|297.13|150.0|10086.3097218415
3d13f1ba-db27-5fa6-9587-1c539c366638|Celso Mourão|COGS|257.9421561647101|404.91|90.0|254.63|297.13|10.0|1314.61215616471
3f4a1eeb-c287-57d9-8407-78ea76ba96ea|Cristina Meireles|R&D|6183.77422031294|364.76|90.0|0.0|297.13|10.0|6945.66422031294
41cf3515-ef53-595d-a1e6-4d127d75533b|Adriana Nogueira|R&D|3719.103451680459|855.8|90.0|0.0|297.13|150.0|5112.033451680459
442868a8-fabc-5568-a849-51f72b08b3c7|Amanda Pinto|R&D|4021.833629163991|848.05|90.0|0.0|297.13|150.0|5407.013629163991
4b86c62c-6ddf-539d-ac69-cc6cc522386a|Diego Figueiredo|G&A|490.174359132607|622.41|90.0|254.63|297.13||1754.344359132607
52f4cadd-4009-5089-b7fa-65f6e4ab2887|Carlos Magalhães|R&D|8837.775201541868|364.76|90.0|254.63|297.13|10.0|9854.295201541867
56f44ead-fafa-5e60-a0a7-e51d81a79de8|Amélia Ribeiro|R&D|5244.623025188975|533.99|117.0|0.0|297.13|150.0|6342.743025188975
607d85fb-c6b1-57ad-8bbe-b1a93fa57349|Benedito Silveira|R&D|1362.231969577253|444.99|90.0|0.0|297.13|150.0|2344.351969577253
6156d6e1-5c40-5ed7-9144-6856fbfbb067|Andreia Peixoto|R&D|8029.53768092652|444.99|90.0|141.46|297.13|150.0|9153.117680926518
61a920cf-a29d-5536-a8ab-3f23dcee4a4b|Eliane Jurema|COGS|4055.50889721234|533.99|90.0|254.63|297.13|50.0|5281.25889721234
631aa093-f502-5229-8a89-b7392dcec06c|André da Rosa|R&D|7487.999892436921|533.99|90.0|141.46|297.13|150.0|8700.57989243692
64bccc6b-6f6b-5e1b-865c-149190776323|Caio Arruda|COGS|8187.839600815882|364.76|246.0|254.63|297.13||9350.35960081588
674419e0-a6fb-5220-acf6-f11f76e231c3|Bernardo Frota|R&D|6387.293051889353|364.76|90.0|0.0|297.13|150.0|7289.183051889353
69b08616-ccd0-5227-b382-8bb0265f456f|Cláudio Dutra|COGS|3997.081814681696|444.99|90.0|0.0|297.13|10.0|4839.201814681695
71019c83-dffc-5bad-8c7f-bd61b699add5|Arthur Gusmão|COGS|5288.260196934793|444.99|90.0|0.0|297.13|150.0|6270.380196934793
7339364e-8890-571d-86f0-ddf62dd4b918|Douglas Sardinha|COGS|2741.365342728207|444.99|90.0|0.0|297.13|50.0|3623.485342728207
73f85436-2fc0-5fc4-840f-f9f03c4fc92b|Aline da Rocha|R&D|1315.118209620667|444.99|90.0|254.63|297.13|150.0|2551.868209620667
7b411419-a6bd-5dac-9123-81532f74e720|Elaine Assis|R&D|5960.861253655476|533.99|90.0|254.63|297.13|50.0|7186.611253655476
95042a82-a9a8-552e-98e3-ac15d20bbf2b|Diogo Quintela|R&D|1461.243910515265|444.99|90.0|0.0|297.13|50.0|2343.363910515265
9ae38b1f-abda-5532-a959-7bfada0e0783|Daniela Drummond|COGS|7065.151216079629|684.63|90.0|0.0|297.13|10.0|8146.911216079629
a0b91621-d2a8-537b-95c4-7ade4c6200ee|Edson Furtado|R&D|7336.902571199591|364.76|90.0|0.0|297.13|50.0|8138.792571199591
a4b3dbb1-86f1-5d25-95a4-82a6c3076da1|Alberto da Luz|R&D|7364.526801321949|444.99|90.0|0.0|297.13|150.0|8346.64680132195
acd778d2-b543-51f4-953c-f9e7e12e7437|Benjamin Queiroz|R&D|7920.453739691342|364.76|206.57|141.46|297.13|150.0|9080.37373969134
b1f0b647-06ea-5aaf-9ea0-a84b2b92563f|Edite Marinho|COGS|800.1804264146962|364.76|90.0|254.63|297.13|50.0|1856.700426414696
b298aa3c-d71a-5df0-83c6-312882b5306b|Cássia Tavares|R&D|8497.84566809797|364.76|90.0|141.46|297.13|10.0|9401.19566809797
b2c05943-feee-5add-9b8e-6bdded9a80a9|Beatriz Siqueira|R&D|3230.18657849905|533.99|90.0|0.0|297.13|150.0|4301.30657849905
b60d0f65-4360-5bc9-b58b-4056c9de1119|Bianca Saraiva|COGS|6015.101383248207|364.76|90.0|254.63|297.13|150.0|7171.621383248207
b79ce13c-6117-58b6-bc9e-06eaab0a22bb|Elias Henriques|R&D|7306.398876912478|-17.8|90.0|0.0|297.13|50.0|7725.728876912478
b96b9339-e1b9-5d20-b70d-67a8bdb16875|Alan Esteves|R&D|8059.789738823023|444.99|90.0|0.0|297.13|150.0|9041.909738823022
b9f4ce97-2bdc-5c19-8c35-5db4419b4948|Bruna Padilha|COGS|7266.55679736717|533.99|90.0|254.63|297.13||8442.30679736717
c4998a9b-cc38-5851-9223-7cf43ea91c36|Alice Farias|G&A|8284.136836853117|560.7|90.0|254.63|297.13|150.0|9636.596836853116
c4fd9f54-9504-52b3-9fca-ab39c50dc07f|Anselmo Novaes|R&D|1108.340345861018|444.99|90.0|0.0|297.13|150.0|2090.460345861018
c57f15e1-0b4a-5ab2-88c8-246678ff4adf|César Franco|R&D|5902.434171124831|0.0|90.0|0.0|297.13|10.0|6299.564171124831
c8fdc760-5724-5f8c-82ae-2e6e767b3811|Denis Macedo|G&A|3952.204962932551|364.76|164.0|254.63|297.13||5032.724962932551
ca2f5a2d-c388-5837-ae1f-93173331ba9a|Alexandre da Mata|R&D|9551.21905794158|364.76|90.0|254.63|297.13|150.0|10707.73905794158
ca78df0b-cb88-5905-ab6f-c47150af2c20|Débora Guedes|COGS|8733.661290965678|855.8|90.0|254.63|297.13|10.0|10241.22129096568
ce6652c4-e8d1-5805-81c2-2fdc603f51b0|Antônio Viana|R&D|3295.032751492262|560.7|117.0|254.63|297.13|150.0|4674.492751492262
d3819e03-5293-5a6a-82ec-d09fd49e456a|Clarice Barreto|R&D|9925.179313508823|364.76|732.0|254.63|297.13|10.0|11583.69931350882
d4d9f427-c9cd-5935-a88d-60bb469e7c89|Álvaro da Paz|R&D|1223.312552727439|533.99|90.0|254.63|297.13|150.0|2549.062552727439
d5426380-369f-500b-a136-ce5543c13c0d|Edgar Cordeiro|R&D|4455.059040524839|444.99|90.0|0.0|297.13|50.0|5337.179040524838
da7b82a2-d4bc-5993-8e5f-59aa3081241a|Daniel Bezerra|R&D|6173.403303137569|444.99|90.0|0.0|297.13|10.0|7015.523303137569
db0b9e18-aa43-5a04-b686-ac454728ad3a|Elielson Quadros|R&D|6715.641722399166|498.39|90.0|254.63|297.13||7855.791722399166
ea3147dc-902e-5ded-88e8-10514c53e142|Ana Melo|R&D|8057.257421621038|444.99|90.0|0.0|297.13|150.0|9039.377421621037
ec7d2162-5de8-50e9-9aea-a1dd90ec8081|Davi Caldeira|R&D|746.9510783679258|560.7|90.0|254.63|297.13|10.0|1959.411078367926
f08290a7-3372-5637-9244-a80f8bc9ee3c|Ângela Rezende|COGS|2036.514997025348|364.76|90.0|0.0|297.13|150.0|2938.404997025348
f35e17a9-38ba-5aa7-9a45-c67338c73cc5|Célia Lemos|R&D|1212.167158944151|364.76|90.0|0.0|297.13|10.0|1974.057158944151
f98cdcb6-f68d-5086-b480-7490f8a1675f|Cristiano Vasconcelos|R&D|5244.95395791246|533.99|90.0|0.0|297.13|10.0|6176.07395791246
f9e02392-bc91-5868-91d5-4a90ccc2cbbf|Anderson Aragão|G&A|9712.728285382607|533.99|90.0|141.46|297.13|150.0|10925.30828538261
"@

$dataLines = $dataBlock -split "`r?`n" | Where-Object { $_.Length -gt 0 }

$rowCount = $dataLines.Count
$arr = New-Object 'object[,]' $rowCount,10

for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $dataLines[$i] -split '\|', -1

    $uuid = $parts[0]
    $name = $parts[1]
    $dept = $parts[2]

    $arr[$i,0] = $uuid
    $arr[$i,1] = $name
    $arr[$i,2] = $dept
    for ($c = 3; $c -lt 10; $c++) {
        $txt = $parts[$c]
        if ($txt -eq "") {
            $arr[$i,$c] = ""
        } else {
            $arr[$i,$c] = [double]$txt
        }
    }
}

$ws.Range("B2:K63").Value = $arr

# Column A keeps the running index 0..61 (unchanged row-for-row; only
# the row CONTENTS were re-sorted by id_unico, the index column is
# positional).
$idxArr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $idxArr[$i,0] = $i
}
$ws.Range("A2:A63").Value = $idxArr

# ---------------------------------------------------------------------
# 3. New trailing row 64: flat "Custo Claude" team-plan line with no
#    colaborador attached (id/nome/depto/salario/outros custos blank).
# ---------------------------------------------------------------------
$ws.Range("A64").Value = 62
$ws.Range("J64").Value = 50
$ws.Range("K64").Value = 50

# Match row 64's index cell (A64) to the same bold/border/centered style
# used by the rest of column A.
$ws.Range("A63").Copy() | Out-Null
$ws.Range("A64").PasteSpecial(-4122) | Out-Null
$ws.Range("A64").Value = 62
